$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Diff shows cell C10 changing from the integer 18 to the numeric value 1
# (serialized as 1.0). Update the cell value accordingly.
$ws.Range("C10").Value = 1
